$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 199.5
$ws.Range("I18").Value = 199.5
$ws.Range("K18").Value = 199.5
$ws.Range("M18").Value = 84.5
$ws.Range("H137").Value = 1225.6757
$ws.Range("I137").Value = 1065.5714
$ws.Range("J137").Value = 1435.8125
$ws.Range("K137").Value = 3196.7142
$ws.Range("L137").Value = 4307.4375
$ws.Range("M137").Value = -646.7142000000003
$ws.Range("N137").Value = -9407.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6945816
$ws.Range("I61").Value = 6945816
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6945816
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -6945604
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 2101
$ws.Range("I74").Value = 2135.1667
$ws.Range("K74").Value = 2135.1667
$ws.Range("M74").Value = -1261.1667
$ws.Range("H77").Value = 2101
$ws.Range("I77").Value = 2135.1667
$ws.Range("K77").Value = 10675.8335
$ws.Range("M77").Value = -6307.833500000001
$ws.Range("H124").Value = 23464.5
$ws.Range("J124").Value = 23464.5
$ws.Range("L124").Value = 23464.5
$ws.Range("N124").Value = -33284.5
$ws.Range("H125").Value = 54480
$ws.Range("J125").Value = 54480
$ws.Range("L125").Value = 54480
$ws.Range("N125").Value = -64320
$ws.Range("H132").Value = 1783672.2
$ws.Range("I132").Value = 1053.4839
$ws.Range("J132").Value = 29414264
$ws.Range("K132").Value = 3160.4517
$ws.Range("L132").Value = 88242792
$ws.Range("M132").Value = -630.4516999999996
$ws.Range("N132").Value = -88247852
$ws.Range("H136").Value = 6945816
$ws.Range("I136").Value = 6945816
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 20837448
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -20834898
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 18519178
$ws.Range("I107").Value = 21739726
$ws.Range("J107").Value = 1022.5
$ws.Range("K107").Value = 21739726
$ws.Range("L107").Value = 1022.5
$ws.Range("M107").Value = -21737806
$ws.Range("N107").Value = -4862.5
$ws.Range("H134").Value = 2821.9092
$ws.Range("I134").Value = 958.53845
$ws.Range("J134").Value = 7363.875
$ws.Range("K134").Value = 2875.61535
$ws.Range("L134").Value = 22091.625
$ws.Range("M134").Value = -340.61535
$ws.Range("N134").Value = -27161.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1185378.4
$ws.Range("I31").Value = 1920057.5
$ws.Range("K31").Value = 1920057.5
$ws.Range("M31").Value = -1919762.5
$ws.Range("H34").Value = 1185378.4
$ws.Range("I34").Value = 1920057.5
$ws.Range("K34").Value = 1920057.5
$ws.Range("M34").Value = -1919855.5
$ws.Range("H48").Value = 7017
$ws.Range("J48").Value = 7017
$ws.Range("L48").Value = 7017
$ws.Range("N48").Value = -7969
$ws.Range("H58").Value = 35715050
$ws.Range("I58").Value = 41667360
$ws.Range("K58").Value = 41667360
$ws.Range("M58").Value = -41667157
$ws.Range("H80").Value = 20000
$ws.Range("J80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("N80").Value = -22246
$ws.Range("H83").Value = 20000
$ws.Range("J83").Value = 20000
$ws.Range("L83").Value = 60000
$ws.Range("N83").Value = -71232
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 13890939
$ws.Range("I132").Value = 1346.7333
$ws.Range("J132").Value = 37040260
$ws.Range("K132").Value = 4040.199900000001
$ws.Range("L132").Value = 111120780
$ws.Range("M132").Value = -1510.199900000001
$ws.Range("N132").Value = -111125840
$ws.Range("H134").Value = 1393.4231
$ws.Range("I134").Value = 1125.1428
$ws.Range("J134").Value = 2520.2
$ws.Range("K134").Value = 3375.4284
$ws.Range("L134").Value = 7560.599999999999
$ws.Range("M134").Value = -840.4284000000002
$ws.Range("N134").Value = -12630.6
$ws.Range("H136").Value = 35715050
$ws.Range("I136").Value = 41667360
$ws.Range("K136").Value = 125002080
$ws.Range("M136").Value = -124999530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 803.9400000000001
$ws.Range("J131").Value = 822.80853
$ws.Range("L131").Value = 2468.42559
$ws.Range("N131").Value = -12548.42559

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 37100
$ws.Range("J15").Value = 37100
$ws.Range("L15").Value = 37100
$ws.Range("N15").Value = -37676
$ws.Range("H81").Value = 37100
$ws.Range("J81").Value = 37100
$ws.Range("L81").Value = 37100
$ws.Range("N81").Value = -39096
$ws.Range("H84").Value = 37100
$ws.Range("J84").Value = 37100
$ws.Range("L84").Value = 111300
$ws.Range("N84").Value = -121284
$ws.Range("H132").Value = 5295.737
$ws.Range("I132").Value = 2418.2856
$ws.Range("J132").Value = 13352.6
$ws.Range("K132").Value = 7254.8568
$ws.Range("L132").Value = 40057.8
$ws.Range("M132").Value = -4724.8568
$ws.Range("N132").Value = -45117.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1474707.2
$ws.Range("I132").Value = 2166115.2
$ws.Range("J132").Value = 22750.4
$ws.Range("K132").Value = 6498345.600000001
$ws.Range("L132").Value = 68251.20000000001
$ws.Range("M132").Value = -6495815.600000001
$ws.Range("N132").Value = -73311.20000000001
$ws.Range("H136").Value = 35715936
$ws.Range("J136").Value = 142858300
$ws.Range("L136").Value = 428574900
$ws.Range("N136").Value = -428580000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 24979.857
$ws.Range("J75").Value = 24979.857
$ws.Range("L75").Value = 24979.857
$ws.Range("N75").Value = -26851.857
$ws.Range("H78").Value = 24979.857
$ws.Range("J78").Value = 24979.857
$ws.Range("L78").Value = 74939.571
$ws.Range("N78").Value = -84299.571
$ws.Range("H86").Value = 12000
$ws.Range("J86").Value = 12000
$ws.Range("L86").Value = 12000
$ws.Range("N86").Value = -14246
$ws.Range("H89").Value = 12000
$ws.Range("J89").Value = 12000
$ws.Range("L89").Value = 60000
$ws.Range("N89").Value = -71232
$ws.Range("H132").Value = 36391.605
$ws.Range("I132").Value = 41439.383
$ws.Range("J132").Value = 17642.715
$ws.Range("K132").Value = 124318.149
$ws.Range("L132").Value = 52928.145
$ws.Range("M132").Value = -121788.149
$ws.Range("N132").Value = -57988.145
$ws.Range("H136").Value = 20002124
$ws.Range("I136").Value = 26317084
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 78951252
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = -78948702
$ws.Range("N136").Value = -19350
